$d = $word.ActiveDocument

# 1) " A No. {1}" -> " A No. [NUMERO_DE_SRA]"
$d.Content.Find.Execute("{1}", $false, $false, $false, $false, $false, $true, 1, $false, "[NUMERO_DE_SRA]", 2)

# 2) "Nombre: {2}" -> "Nombre: [NOMBRE_CLIENTE]"
$d.Content.Find.Execute("{2}", $false, $false, $false, $false, $false, $true, 1, $false, "[NOMBRE_CLIENTE]", 2)

# 3) "N de Servicio: {3}" -> "N de Servicio: [NUMERO_DE_SR]"
$d.Content.Find.Execute("{3}", $false, $false, $false, $false, $false, $true, 1, $false, "[NUMERO_DE_SR]", 2)

# 4) Move the month placeholder: "Zona No. {4} {5}<br>En concepto ... mes de Febrero del Año 2023"
#    -> "Zona No. {4}<br>En concepto ... mes de {5} del Año 2023"
#    (there is a manual line break - vertical tab char 11 - between "{5}" and "En concepto")
$br = [char]11
$findText4 = "Zona No. {4} {5}" + $br + "En concepto de Servicio de Agua Potable correspondiente al mes de Febrero del Año 2023"
$replText4 = "Zona No. {4}" + $br + "En concepto de Servicio de Agua Potable correspondiente al mes de {5} del Año 2023"
$d.Content.Find.Execute($findText4, $false, $false, $false, $false, $false, $true, 1, $false, $replText4, 2)

# 5) Cuota Normal de Agua potable amount: "{6}" -> "x"
$d.Content.Find.Execute("{6}", $false, $false, $false, $false, $false, $true, 1, $false, "x", 2)

# 6) Cuota por Agregados amount: "{7}" -> "x"
$d.Content.Find.Execute("{7}", $false, $false, $false, $false, $false, $true, 1, $false, "x", 2)

# 7) cancelación de Multas amount: "{8}" -> "x"
$d.Content.Find.Execute("{8}", $false, $false, $false, $false, $false, $true, 1, $false, "x", 2)

# 8) Colaboraciones amount: "{9}" -> "x"
$d.Content.Find.Execute("{9}", $false, $false, $false, $false, $false, $true, 1, $false, "x", 2)

# 9) TOTAL A PAGAR amount: "{10}" -> "x"  (spans two runs in the source, Find still matches across them)
$d.Content.Find.Execute("{10}", $false, $false, $false, $false, $false, $true, 1, $false, "x", 2)

# 10) "Fecha Limite de Pago: {12}" -> "Fecha Limite de Pago: 2023-06-30"
$d.Content.Find.Execute("{12}", $false, $false, $false, $false, $false, $true, 1, $false, "2023-06-30", 2)
